# Insert a new weekly data row at row 14 (pushes all subsequent rows down by
# one), matching the commit "Fruta / hortaliza, semanal" which appends a new
# sampling date ahead of the previously-newest rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("14:14").Insert()

$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44558
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112040
$ws.Range("G14").Value = "Cilantro"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 3200
$ws.Range("K14").Value = 2300
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2400
$ws.Range("N14").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 1600
$ws.Range("Q14").Value = 1.5
$ws.Range("R14").Value = "Hortaliza"
